# Updates the cryptocurrency price/volume table (and a couple of coin
# name/link swaps) on the active worksheet to match the refreshed
# "cryptos list" snapshot, per the GitHub Actions data-refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new text value, and whether the cell must be
# forced to Text format first so Excel does not silently reinterpret a
# numeric-looking string (e.g. "299.21") as a Number, which would change
# both its stored type and drop things like trailing zeros (e.g. "2.80").
$updates = @(
    @{ Cell = 'D2'; Value = '42.111.46'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  -1.32%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '2.268.91'; ForceText = $false }
    @{ Cell = 'E3'; Value = '  -1.56%  '; ForceText = $false }
    @{ Cell = 'E4'; Value = '  +0.08%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '299.21'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  -1.36%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '95.41'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  -4.22%  '; ForceText = $false }
    @{ Cell = 'D7'; Value = '0.494'; ForceText = $true }
    @{ Cell = 'E7'; Value = '  -2.51%  '; ForceText = $false }
    @{ Cell = 'E8'; Value = '  +0.05%  '; ForceText = $false }
    @{ Cell = 'E9'; Value = '  -2.82%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '33.13'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  -4.74%  '; ForceText = $false }
    @{ Cell = 'D11'; Value = '0.0787'; ForceText = $true }
    @{ Cell = 'E11'; Value = '  -0.67%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '48.48'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  -6.51%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '0.114'; ForceText = $true }
    @{ Cell = 'E13'; Value = '  +0.66%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '6.67'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  -1.62%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '15.77'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  +0.52%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '2.622.99'; ForceText = $false }
    @{ Cell = 'E16'; Value = '  -1.43%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '2.264.41'; ForceText = $false }
    @{ Cell = 'E17'; Value = '  -1.89%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '0.783'; ForceText = $true }
    @{ Cell = 'E18'; Value = '  -2.50%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '42.102.56'; ForceText = $false }
    @{ Cell = 'E19'; Value = '  -1.18%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '11.76'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  +2.19%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '0.0₃0891'; ForceText = $false }
    @{ Cell = 'E21'; Value = '  -1.39%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '5.97'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  -1.24%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '66.12'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  -3.03%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '235.16'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  -0.21%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '1.95'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  -0.86%  '; ForceText = $false }
    @{ Cell = 'E26'; Value = '  +0.12%  '; ForceText = $false }
    @{ Cell = 'E27'; Value = '  -2.48%  '; ForceText = $false }
    @{ Cell = 'D28'; Value = '23.71'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  -5.84%  '; ForceText = $false }
    @{ Cell = 'D29'; Value = '2.16'; ForceText = $true }
    @{ Cell = 'E29'; Value = '  +4.37%  '; ForceText = $false }
    @{ Cell = 'D30'; Value = '168.22'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  +3.26%  '; ForceText = $false }
    @{ Cell = 'D31'; Value = '9.17'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  -0.35%  '; ForceText = $false }
    @{ Cell = 'D32'; Value = '33.68'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  -3.29%  '; ForceText = $false }
    @{ Cell = 'D33'; Value = '0.999'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  -0.04%  '; ForceText = $false }
    @{ Cell = 'D34'; Value = '4.88'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  -2.59%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '4.57'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  -0.83%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '16.71'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  -1.55%  '; ForceText = $false }
    @{ Cell = 'E37'; Value = '  -2.80%  '; ForceText = $false }
    @{ Cell = 'D38'; Value = '0.0687'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  -3.51%  '; ForceText = $false }
    @{ Cell = 'D39'; Value = '2.80'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  -2.91%  '; ForceText = $false }
    @{ Cell = 'E40'; Value = '  -1.56%  '; ForceText = $false }
    @{ Cell = 'B41'; Value = 'ARBITRUM'; ForceText = $false }
    @{ Cell = 'C41'; Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; ForceText = $false }
    @{ Cell = 'D41'; Value = '1.72'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  -5.00%  '; ForceText = $false }
    @{ Cell = 'B42'; Value = 'Stellar'; ForceText = $false }
    @{ Cell = 'C42'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; ForceText = $false }
    @{ Cell = 'D42'; Value = '0.108'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  -2.93%  '; ForceText = $false }
    @{ Cell = 'D43'; Value = '2.34'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  -3.79%  '; ForceText = $false }
    @{ Cell = 'D44'; Value = '1.963.95'; ForceText = $false }
    @{ Cell = 'E44'; Value = '  -0.49%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '0.0277'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  -1.27%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '17.56'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  -4.96%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '9.55'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  -6.60%  '; ForceText = $false }
    @{ Cell = 'D48'; Value = '2.77'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  -4.54%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '2.494.50'; ForceText = $false }
    @{ Cell = 'E49'; Value = '  -1.31%  '; ForceText = $false }
    @{ Cell = 'D50'; Value = '52.17'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  -5.43%  '; ForceText = $false }
    @{ Cell = 'B51'; Value = 'Stacks'; ForceText = $false }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'; ForceText = $false }
    @{ Cell = 'D51'; Value = '1.48'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  -0.76%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $rng.NumberFormat = '@'
    }
    $rng.Value = $u.Value
}
